# "Doing Updates for Financials"
# Update the D column (most recent fiscal period, 2018-06-30) figures on the
# TRT sheet: several cells that previously held the placeholder "NA" text now
# have real numeric values, and a handful of existing numbers were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income statement section
$ws.Range("D17").Value = 40200    # Total Operating Expenses
$ws.Range("D18").Value = 2200     # Operating Income or Loss
$ws.Range("D20").Value = 300      # Total Other Income/Expenses Net
$ws.Range("D21").Value = 4700     # Earnings Before Interest And Taxes
$ws.Range("D24").Value = 1000     # Income Tax Expense
$ws.Range("D26").Value = 1300     # Income After Tax
$ws.Range("D27").Value = 1200     # Net Income From Continuing Ops
$ws.Range("D32").Value = -300     # Other Items
$ws.Range("D33").Value = 1200     # Net Income
$ws.Range("D35").Value = 1200     # Net Income Applicable To Common Shares

# Balance sheet section
$ws.Range("D59").Value = 3500     # Other Current Liabilities (was 4300)
$ws.Range("D62").Value = 1200     # Other Liabilities (was 400)
$ws.Range("D76").Value = 22000    # Total Stockholder Equity (was 23500)

# Cash flow statement section
$ws.Range("D81").Value = 1200     # Net Income
